# Analyse de la répartition du temps de travail II
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3) ---
$ws.Range("E3").Value = "Adrien"
$ws.Range("F3").Value = "James"
$ws.Range("D3").Value = "Total"
$ws.Range("G3").Value = "Antoine"
$ws.Range("H3").Value = "Loyse"
$ws.Range("I3").Value = "Jérémie"

# --- Row 4 : Planification ---
$ws.Range("D4").Formula = "=SUM(E4:J4)"
$ws.Range("E4").Formula = "=3+1+1+5+1+1.5+2+2+1"
$ws.Range("F4").Formula = "=11+1.5+1.5"
$ws.Range("G4").Formula = "=10.5+3.5+1.5+1.5+1.5"
$ws.Range("H4").Formula = "=4.75+3+1.5"
$ws.Range("I4").Formula = "=12.5+3+3"

# --- Row 5 : Formation ---
$ws.Range("D5").Formula = "=SUM(E5:J5)"
$ws.Range("E5").Formula = "=5+2"
$ws.Range("F5").Formula = "=4+4"
$ws.Range("G5").Formula = "=6+3"
$ws.Range("H5").Formula = "=9+1.5+1.5"
$ws.Range("I5").Formula = "=10+3"

# --- Row 6 : Modélisation ---
$ws.Range("D6").Formula = "=SUM(E6:J6)"
$ws.Range("E6").Formula = "=5+2.5+4"
$ws.Range("F6").Value = 0
$ws.Range("G6").Formula = "=8+4"
$ws.Range("H6").Formula = "=12+8"
$ws.Range("I6").Formula = "=6+4"

# --- Row 7 : GUI ---
$ws.Range("D7").Formula = "=SUM(E7:J7)"
$ws.Range("E7").Value = 0
$ws.Range("F7").Formula = "=6+5+5.5+1.5"
$ws.Range("G7").Formula = "=4+2+3"
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# --- Row 8 : Sauvegarde ---
$ws.Range("D8").Formula = "=SUM(E8:J8)"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Formula = "=18.5+5"

# --- Row 9 : Historique ---
$ws.Range("D9").Formula = "=SUM(E9:J9)"
$ws.Range("E9").Formula = "=5+ 13+10"
$ws.Range("F9").Value = 8
$ws.Range("G9").Formula = "=3"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 24

# --- Row 10 : Outils ---
$ws.Range("D10").Formula = "=SUM(E10:J10)"
$ws.Range("E10").Formula = "=4+2+2"
$ws.Range("F10").Formula = "=0.5+36.5+1"
$ws.Range("G10").Formula = "=5+1+12+10+12+3+7"
$ws.Range("H10").Formula = "=9+13+10"
$ws.Range("I10").Value = 0

# --- Row 11 : Documentation ---
$ws.Range("D11").Formula = "=SUM(E11:J11)"
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 5.5
$ws.Range("G11").Formula = "=3+6"
$ws.Range("H11").Formula = "=16"
$ws.Range("I11").Formula = "=21"

# --- Row 12 : Autres ---
$ws.Range("D12").Formula = "=SUM(E12:J12)"
$ws.Range("E12").Value = 7
$ws.Range("F12").Formula = "=0.5+4.5"
$ws.Range("G12").Formula = "=2"
$ws.Range("H12").Value = 0.25
$ws.Range("I12").Value = 0

# --- Row 13 : Présentation ---
$ws.Range("D13").Formula = "=SUM(E13:J13)"
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# --- Row 15 : total ---
$ws.Range("F15").Formula = "=SUM(F4:F13)"
$ws.Range("G15").Formula = "=SUM(G4:G13)"
$ws.Range("H15").Formula = "=SUM(H4:H13)"
$ws.Range("I15").Formula = "=SUM(I4:I13)"
$ws.Range("J15").Formula = "=SUM(J4:J13)"

# --- Workbook calculation options (enable iterative calc, tighten max change) ---
$excel.IterativeCalculation = $true
$excel.MaxChange = 0.0001

# --- View tweaks ---
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I19").Select()

$wb.Save()
